# Updated via Streamlit Approval System
# Adds ACCEPTED / PAID / HOLD / REJECTED checkbox columns (AP:AS) to the
# pending-approval sheet, flags rows 5-16 as HOLD, and clears the legacy
# text-based HOLD markers in row 5 (APPROVAL_1 / APPROVAL_2 columns).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells (row 1), matching the style of the existing headers ---
$ws.Range("AO1").Copy()
$ws.Range("AP1:AS1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("AP1").Value = "ACCEPTED"
$ws.Range("AQ1").Value = "PAID"
$ws.Range("AR1").Value = "HOLD"
$ws.Range("AS1").Value = "REJECTED"

# --- New boolean checkbox columns for every data row ---
$lastRow = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 42).Value = $false   # AP - ACCEPTED
    $ws.Cells.Item($r, 43).Value = $false   # AQ - PAID
    $ws.Cells.Item($r, 45).Value = $false   # AS - REJECTED
}

# Rows 5-16 are on hold
for ($r = 5; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 44).Value = $true    # AR - HOLD
}
# Rows 2-4 are not on hold
for ($r = 2; $r -le 4; $r++) {
    $ws.Cells.Item($r, 44).Value = $false   # AR - HOLD
}

# --- Clear the legacy text-based HOLD markers for row 5 now that the ---
# --- checkbox-based HOLD flag (AR5) captures the state ---
$ws.Range("AI5").Value = ""
$ws.Range("AJ5").Value = ""
